$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that currently sits near the end
#    of the document (inside the last table cell, right after the word
#    "numerical").
# ---------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# ---------------------------------------------------------------------
# 2) Remove the curly quotes that wrap "Socio Demographic Indicators"
#    in the very first paragraph, turning:
#       "Socio Demographic Indicators":
#    into:
#       Socio Demographic Indicators:
#    The opening quote is the 1st character (offset 0-1); the closing
#    quote sits right before the final colon. Delete the closing quote
#    first so the opening quote's offsets stay valid.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$closeQuoteStart = $p1.Text.IndexOf([char]8221)
$closeQuote = $d.Range($p1.Start + $closeQuoteStart, $p1.Start + $closeQuoteStart + 1)
$closeQuote.Delete()

$openQuote = $d.Range($p1.Start, $p1.Start + 1)
$openQuote.Delete()

# ---------------------------------------------------------------------
# 3) Insert a new "_GoBack" bookmark collapsed at the very start of the
#    document (before the word "Socio"). The engine has trouble
#    collapsing a bookmark exactly at offset 0 directly, so use a
#    small, reliable workaround: insert a temporary marker character at
#    the start, add the bookmark right after it, then remove the
#    marker.
# ---------------------------------------------------------------------
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("@")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$markerRange = $d.Range(0, 1)
$markerRange.Delete()

# ---------------------------------------------------------------------
# 4) Color some analysis-variable cells red (FF0000) in tables 1 and 2.
# ---------------------------------------------------------------------

# Table 1, row 4: "Rural (% pop in 2010), County Health Ranking 2019" /
#                 "????????????????????"
$tbl1 = $d.Tables.Item(1)
$tbl1.Cell(4, 1).Range.Font.Color = 255
$tbl1.Cell(4, 2).Range.Font.Color = 255

# Table 2, row 4: "Median Income SAIPE ..." / "Median Income ($ in 2017) ..."
$tbl2 = $d.Tables.Item(2)
$tbl2.Cell(4, 1).Range.Font.Color = 255
$tbl2.Cell(4, 2).Range.Font.Color = 255

# ---------------------------------------------------------------------
# 5) Delete the row "NEVERMIND NOD CNTY DATA .....SNAP participants
#    (% pop in 2016), USDA 2017" / "?????????????" entirely (table 3,
#    row 3).
# ---------------------------------------------------------------------
$tbl3 = $d.Tables.Item(3)
$tbl3.Rows.Item(3).Delete()
